# Regression test for configuration file validation:
# add a new "solidity" column to the "geometry" sheet, inserted right
# before the existing "thickness_max_chord_ratio" column (AN), shifting
# that and all following columns one to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Insert a new column at AN; this shifts AN:AQ -> AO:AR and keeps the
# header style (bold font + border) that was on the old AN column.
$ws.Range("AN1:AN2").EntireColumn.Insert()

# Populate the newly inserted column.
$ws.Range("AN1").Value = "solidity"
$ws.Range("AN2").Value = "[1.42997704 1.70997375]"
